# Fix closing dates / rate values for algo backup sheet (BAJFIN)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 7 (daily summary row)
$ws.Range("F7").Value = 7445.3
$ws.Range("G7").Value = 7468.95
$ws.Range("H7").Value = 7305.5
$ws.Range("I7").Value = 7367.15
$ws.Range("J7").Value = 7493.95

# Row 9
$ws.Range("G9").Value = 7477
$ws.Range("H9").Value = 7407
$ws.Range("I9").Value = 7421

# Row 10
$ws.Range("G10").Value = 7439.75
$ws.Range("H10").Value = 7401
$ws.Range("I10").Value = 7428

# Row 11
$ws.Range("G11").Value = 7461.95
$ws.Range("H11").Value = 7419.1
$ws.Range("I11").Value = 7432.15

# Row 12
$ws.Range("G12").Value = 7449.95
$ws.Range("H12").Value = 7415
$ws.Range("I12").Value = 7443.9

# Row 13
$ws.Range("G13").Value = 7457.45
$ws.Range("H13").Value = 7426
$ws.Range("I13").Value = 7441.35

# Row 14
$ws.Range("G14").Value = 7453.85
$ws.Range("H14").Value = 7430.3
$ws.Range("I14").Value = 7433.5

# Row 15
$ws.Range("G15").Value = 7448.9
$ws.Range("H15").Value = 7428.1
$ws.Range("I15").Value = 7439.65

# Row 16
$ws.Range("G16").Value = 7468.95
$ws.Range("H16").Value = 7433.9
$ws.Range("I16").Value = 7459.15

# Row 17
$ws.Range("G17").Value = 7468.05
$ws.Range("H17").Value = 7415.5
$ws.Range("I17").Value = 7420.05

# Row 18
$ws.Range("G18").Value = 7417.35
$ws.Range("H18").Value = 7350
$ws.Range("I18").Value = 7362.2

# Row 19
$ws.Range("G19").Value = 7379.2
$ws.Range("H19").Value = 7305.5
$ws.Range("I19").Value = 7353.95

# Row 20
$ws.Range("G20").Value = 7367.55
$ws.Range("H20").Value = 7342.2
$ws.Range("I20").Value = 7349.7

# Row 21
$ws.Range("G21").Value = 7388
$ws.Range("H21").Value = 7346.1
$ws.Range("I21").Value = 7386
